$d = $word.ActiveDocument

# 1. Vet paragraph: insert " и управляет процессами размещения животных на передержку" after "передержке"
$d.Content.Find.Execute(
    "Ветеринар осуществляет уход за животными на передержке, проводит регулярные медицинские осмотры",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Ветеринар осуществляет уход за животными на передержке и управляет процессами размещения животных на передержку, проводит регулярные медицинские осмотры",
    2
) | Out-Null

# 2. Manager paragraph: replace tail
$d.Content.Find.Execute(
    "Менеджер управляет процессами размещения животных на передержку, ведет учет поступающих и выбывающих животных.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Менеджер управляет процессами обработки всех платежей платежей.",
    2
) | Out-Null

# 3. Administrator paragraph: re-write as a single run (merges the split runs)
$d.Content.Find.Execute(
    "Администратор управляет системой, назначает права доступа, контролирует корректность ведения учета и поддерживает базу данных.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Администратор управляет системой, назначает права доступа, контролирует корректность ведения учета и поддерживает базу данных.",
    2
) | Out-Null

# 4. Windows/MS Office paragraph: re-write as a single run (merges the split runs, drops proofErr tags)
$d.Content.Find.Execute(
    "Программа должна работать в операционных системах Windows 10/11. Все формируемые отчеты должны иметь возможность экспортирования в редактор электронных таблиц MS Office Excel 2019/2021.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Программа должна работать в операционных системах Windows 10/11. Все формируемые отчеты должны иметь возможность экспортирования в редактор электронных таблиц MS Office Excel 2019/2021.",
    2
) | Out-Null

# 5. Table date insertions (empty cells get new date text)
$tbl = $d.Tables.Item(1)

$tbl.Cell(3,3).Range.Text = "07.03.2025"
$tbl.Cell(3,4).Range.Text = "09.02.2025"

$tbl.Cell(4,3).Range.Text = "10.03.2025"
$tbl.Cell(4,4).Range.Text = "11.03.2025"

$tbl.Cell(5,3).Range.Text = "12.03.2025"
$tbl.Cell(5,4).Range.Text = "13.03.2025"

$tbl.Cell(6,3).Range.Text = "14.03.2025"
$tbl.Cell(6,4).Range.Text = "15.03.2025"

$tbl.Cell(7,3).Range.Text = "16.03.2025"
$tbl.Cell(7,4).Range.Text = "04.04.2025"

$tbl.Cell(8,3).Range.Text = "05.04.2025"
$tbl.Cell(8,4).Range.Text = "28.04.2024"

# 6. Remove the stray _GoBack bookmark
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 7. Defense dates: 25.04.2025 -> 29.04.2025 (twice)
$d.Content.Find.Execute(
    "25.04.2025", $true, $false, $false, $false, $false, $true, 1, $false,
    "29.04.2025", 2
) | Out-Null
$d.Content.Find.Execute(
    "25.04.2025", $true, $false, $false, $false, $false, $true, 1, $false,
    "29.04.2025", 2
) | Out-Null
